$wb = $excel.ActiveWorkbook

# --- Generator sheet: becomes the active tab, H3 value changes 60 -> 100, selection moves to H4 ---
$wsGenerator = $wb.Worksheets.Item("Generator")
$wsGenerator.Range("H3").Value = 100

# --- StorageUnit sheet: L2 value changes 0 -> 20 (selection/tab state stays at L3) ---
$wsStorageUnit = $wb.Worksheets.Item("StorageUnit")
$wsStorageUnit.Range("L2").Value = 20

# Activate the Generator sheet last so it becomes the workbook's active tab
# and select H4 on it, matching the new selection/tabSelected state.
$wsGenerator.Activate()
$wsGenerator.Range("H4").Select()
